# "adapted wireing width of power stage and exchanged bootstrap capacitor"
#
# The power-stage part list gains one more row describing the newly
# exchanged bootstrap capacitor (C1): its schematic designator, its
# description and a Mouser link to the part (formatted/linked the same
# way as the other "Mouser Link" entries in column C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$url = "https://www.mouser.de/ProductDetail/TDK/CGA4J1X7R1H475K125AE?qs=PqoDHHvF64%252BKyTaX2qYTwQ%3D%3D "

# Fill the new row. The link cell is written first, then the description,
# then the designator - this mirrors how the part was originally entered
# and keeps the shared-string table ordering faithful.
$ws.Range("C11").Value = $url
$ws.Range("B11").Value = "C bootstrap"
$ws.Range("A11").Value = "C1"

# Turn the Mouser link into a real hyperlink, same as the other rows in
# column C (this also gives the cell the underlined "link" look used
# throughout the sheet).
[void]$ws.Hyperlinks.Add($ws.Range("C11"), $url)

# Leave the selection where the author ended up after entering the data.
[void]$ws.Range("B14").Select()
